# Natmi following Dr Hou advice
# Update recalculated Il34-Csf1r LR-pair statistics (rows 2-16, columns E-T)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.654587
$ws.Range("H2").Value = 4.963761
$ws.Range("I2").Value = 0.05277312046682094
$ws.Range("J2").Value = 0.07150718413853953
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.654529
$ws.Range("N2").Value = 1.963587
$ws.Range("O2").Value = 0.001636504842097015
$ws.Range("P2").Value = 0.001641975240588762
$ws.Range("Q2").Value = 1.082975174523
$ws.Range("R2").Value = 9.746776570707
$ws.Range("S2").Value = [double]"8.636346717652153E-05"
$ws.Range("T2").Value = 0.0001174130258797034
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.654587
$ws.Range("H3").Value = 4.963761
$ws.Range("I3").Value = 0.05277312046682094
$ws.Range("J3").Value = 0.07150718413853953
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.625878666666667
$ws.Range("N3").Value = 7.877636000000001
$ws.Range("O3").Value = 0.006565428197618827
$ws.Range("P3").Value = 0.006587374670116828
$ws.Range("Q3").Value = 4.344744705444
$ws.Range("R3").Value = 39.102702348996
$ws.Range("S3").Value = 0.0003464781331892014
$ws.Range("T3").Value = 0.0004710446135255951
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.654587
$ws.Range("H4").Value = 4.963761
$ws.Range("I4").Value = 0.05277312046682094
$ws.Range("J4").Value = 0.07150718413853953
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 196.492958
$ws.Range("N4").Value = 589.4788739999999
$ws.Range("O4").Value = 0.491287135031397
$ws.Range("P4").Value = 0.4929293766755139
$ws.Range("Q4").Value = 325.1146938983459
$ws.Range("R4").Value = 2926.032245085114
$ws.Range("S4").Value = 0.02592675516081124
$ws.Range("T4").Value = 0.03524799170523148
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.654587
$ws.Range("H5").Value = 4.963761
$ws.Range("I5").Value = 0.05277312046682094
$ws.Range("J5").Value = 0.07150718413853953
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 196.1845956666667
$ws.Range("N5").Value = 588.5537870000001
$ws.Range("O5").Value = 0.4905161432928793
$ws.Range("P5").Value = 0.4921558077175863
$ws.Range("Q5").Value = 324.6044815903231
$ws.Range("R5").Value = 2921.440334312907
$ws.Range("S5").Value = 0.02588606752091552
$ws.Range("T5").Value = 0.0351926759673131
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.654587
$ws.Range("H6").Value = 4.963761
$ws.Range("I6").Value = 0.05277312046682094
$ws.Range("J6").Value = 0.07150718413853953
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 3.99747
$ws.Range("N6").Value = 7.99494
$ws.Range("O6").Value = 0.009994788636007805
$ws.Range("P6").Value = 0.006685465696194116
$ws.Range("Q6").Value = 6.61416189489
$ws.Range("R6").Value = 39.68497136934
$ws.Range("S6").Value = 0.0005274561847284527
$ws.Range("T6").Value = 0.000478058826589642
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 5.055996666666666
$ws.Range("H7").Value = 15.16799
$ws.Range("I7").Value = 0.1612612217851615
$ws.Range("J7").Value = 0.2185077512679451
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.654529
$ws.Range("N7").Value = 1.963587
$ws.Range("O7").Value = 0.001636504842097015
$ws.Range("P7").Value = 0.001641975240588762
$ws.Range("Q7").Value = 3.309296442236667
$ws.Range("R7").Value = 29.78366798013
$ws.Range("S7").Value = 0.0002639047702938974
$ws.Range("T7").Value = 0.0003587843174586935
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 5.055996666666666
$ws.Range("H8").Value = 15.16799
$ws.Range("I8").Value = 0.1612612217851615
$ws.Range("J8").Value = 0.2185077512679451
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.625878666666667
$ws.Range("N8").Value = 7.877636000000001
$ws.Range("O8").Value = 0.006565428197618827
$ws.Range("P8").Value = 0.006587374670116828
$ws.Range("Q8").Value = 13.27643378573778
$ws.Range("R8").Value = 119.48790407164
$ws.Range("S8").Value = 0.001058748972690763
$ws.Range("T8").Value = 0.001439392425926649
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 5.055996666666666
$ws.Range("H9").Value = 15.16799
$ws.Range("I9").Value = 0.1612612217851615
$ws.Range("J9").Value = 0.2185077512679451
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 196.492958
$ws.Range("N9").Value = 589.4788739999999
$ws.Range("O9").Value = 0.491287135031397
$ws.Range("P9").Value = 0.4929293766755139
$ws.Range("Q9").Value = 993.4677406714731
$ws.Range("R9").Value = 8941.209666043258
$ws.Range("S9").Value = 0.07922556364249472
$ws.Range("T9").Value = 0.1077088896312764
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 5.055996666666666
$ws.Range("H10").Value = 15.16799
$ws.Range("I10").Value = 0.1612612217851615
$ws.Range("J10").Value = 0.2185077512679451
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 196.1845956666667
$ws.Range("N10").Value = 588.5537870000001
$ws.Range("O10").Value = 0.4905161432928793
$ws.Range("P10").Value = 0.4921558077175863
$ws.Range("Q10").Value = 991.9086617420145
$ws.Range("R10").Value = 8927.177955678131
$ws.Range("S10").Value = 0.07910123257275509
$ws.Range("T10").Value = 0.1075398588178289
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 5.055996666666666
$ws.Range("H11").Value = 15.16799
$ws.Range("I11").Value = 0.1612612217851615
$ws.Range("J11").Value = 0.2185077512679451
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 3.99747
$ws.Range("N11").Value = 7.99494
$ws.Range("O11").Value = 0.009994788636007805
$ws.Range("P11").Value = 0.006685465696194116
$ws.Range("Q11").Value = 20.2111949951
$ws.Range("R11").Value = 121.2671699706
$ws.Range("S11").Value = 0.001611771826927067
$ws.Range("T11").Value = 0.001460826075454363
$ws.Range("E12").Value = 2
$ws.Range("G12").Value = 24.6422525
$ws.Range("H12").Value = 49.284505
$ws.Range("I12").Value = 0.7859656577480176
$ws.Range("J12").Value = 0.7099850645935153
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.654529
$ws.Range("N12").Value = 1.963587
$ws.Range("O12").Value = 0.001636504842097015
$ws.Range("P12").Value = 0.001641975240588762
$ws.Range("Q12").Value = 16.1290688865725
$ws.Range("R12").Value = 96.77441331943498
$ws.Range("S12").Value = 0.001286236604626596
$ws.Range("T12").Value = 0.001165777897250365
$ws.Range("E13").Value = 2
$ws.Range("G13").Value = 24.6422525
$ws.Range("H13").Value = 49.284505
$ws.Range("I13").Value = 0.7859656577480176
$ws.Range("J13").Value = 0.7099850645935153
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 2.625878666666667
$ws.Range("N13").Value = 7.877636000000001
$ws.Range("O13").Value = 0.006565428197618827
$ws.Range("P13").Value = 0.006587374670116828
$ws.Range("Q13").Value = 64.70756513836334
$ws.Range("R13").Value = 388.24539083018
$ws.Range("S13").Value = 0.005160201091738863
$ws.Range("T13").Value = 0.004676937630664583
$ws.Range("E14").Value = 2
$ws.Range("G14").Value = 24.6422525
$ws.Range("H14").Value = 49.284505
$ws.Range("I14").Value = 0.7859656577480176
$ws.Range("J14").Value = 0.7099850645935153
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 196.492958
$ws.Range("N14").Value = 589.4788739999999
$ws.Range("O14").Value = 0.491287135031397
$ws.Range("P14").Value = 0.4929293766755139
$ws.Range("Q14").Value = 4842.029085507894
$ws.Range("R14").Value = 29052.17451304736
$ws.Range("S14").Value = 0.3861348162280911
$ws.Range("T14").Value = 0.349972495339006
$ws.Range("E15").Value = 2
$ws.Range("G15").Value = 24.6422525
$ws.Range("H15").Value = 49.284505
$ws.Range("I15").Value = 0.7859656577480176
$ws.Range("J15").Value = 0.7099850645935153
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 196.1845956666667
$ws.Range("N15").Value = 588.5537870000001
$ws.Range("O15").Value = 0.4905161432928793
$ws.Range("P15").Value = 0.4921558077175863
$ws.Range("Q15").Value = 4834.430343028406
$ws.Range("R15").Value = 29006.58205817044
$ws.Range("S15").Value = 0.3855288431992088
$ws.Range("T15").Value = 0.3494232729324442
$ws.Range("E16").Value = 2
$ws.Range("G16").Value = 24.6422525
$ws.Range("H16").Value = 49.284505
$ws.Range("I16").Value = 0.7859656577480176
$ws.Range("J16").Value = 0.7099850645935153
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 3.99747
$ws.Range("N16").Value = 7.99494
$ws.Range("O16").Value = 0.009994788636007805
$ws.Range("P16").Value = 0.006685465696194116
$ws.Range("Q16").Value = 98.50666510117499
$ws.Range("R16").Value = 394.0266604047
$ws.Range("S16").Value = 0.007855560624352286
$ws.Range("T16").Value = 0.00474658079415011
